$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"22.83000000000013"
$ws.Range("H2").Value = [double]"4.271460841698271e-07"
$ws.Range("I2").Value = [double]"4.271460841698271e-07"
$ws.Range("L2").Value = [double]"54.5261417984353"
$ws.Range("M2").Value = "[32.24136979362113, 76.81091380324946]"
$ws.Range("N2").Value = [double]"1.168345493862688e-05"
$ws.Range("O2").Value = [double]"1.168345493862688e-05"
$ws.Range("P2").Value = [double]"2.018921405009658"
$ws.Range("Q2").Value = "[1.603816069400196, 2.434026740619119]"
$ws.Range("R2").Value = [double]"9.883205365213144e-13"
$ws.Range("S2").Value = [double]"9.883205365213144e-13"
$ws.Range("T2").Value = [double]"60.83342505429179"
$ws.Range("U2").Value = "[48.879266225984225, 72.78758388259936]"
$ws.Range("V2").Value = [double]"2.389199948993337e-13"
$ws.Range("W2").Value = [double]"2.389199948993337e-13"
$ws.Range("X2").Value = [double]"15.49423423423432"
$ws.Range("Y2").Value = [double]"13.98594594594602"
$ws.Range("Z2").Value = [double]"17.00252252252261"

$ws.Range("F3").Value = [double]"22.83000000000013"
$ws.Range("H3").Value = [double]"0.01517751521992405"
$ws.Range("I3").Value = [double]"0.01517751521992405"
$ws.Range("L3").Value = [double]"29.24961833122104"
$ws.Range("M3").Value = "[3.5218839266456747, 54.9773527357964]"
$ws.Range("N3").Value = [double]"0.02677395336222266"
$ws.Range("O3").Value = [double]"0.02677395336222266"
$ws.Range("P3").Value = [double]"2.232763547596349"
$ws.Range("Q3").Value = "[1.415131825941347, 3.050395269251351]"
$ws.Range("R3").Value = [double]"1.713112408152639e-06"
$ws.Range("S3").Value = [double]"1.713112408152639e-06"
$ws.Range("T3").Value = [double]"49.81269475255375"
$ws.Range("U3").Value = "[36.088232144508474, 63.53715736059902]"
$ws.Range("V3").Value = [double]"3.541515081195712e-09"
$ws.Range("W3").Value = [double]"3.541515081195712e-09"
$ws.Range("X3").Value = [double]"14.71723723723732"
$ws.Range("Y3").Value = [double]"11.74636636636643"
$ws.Range("Z3").Value = [double]"17.68810810810821"

$ws.Range("F4").Value = [double]"22.83000000000013"
$ws.Range("H4").Value = [double]"0.009850686790391583"
$ws.Range("I4").Value = [double]"0.009850686790391583"
$ws.Range("L4").Value = [double]"28.92562630240473"
$ws.Range("M4").Value = "[3.7196050555180804, 54.13164754929137]"
$ws.Range("N4").Value = [double]"0.02545142985969484"
$ws.Range("O4").Value = [double]"0.02545142985969484"
$ws.Range("P4").Value = [double]"2.647868883205812"
$ws.Range("Q4").Value = "[1.5912371198362725, 3.7045006465753514]"
$ws.Range("R4").Value = [double]"7.858704720309362e-06"
$ws.Range("S4").Value = [double]"7.858704720309362e-06"
$ws.Range("T4").Value = [double]"54.17087607937374"
$ws.Range("U4").Value = "[41.11233268932314, 67.22941946942434]"
$ws.Range("V4").Value = [double]"1.054585307969091e-10"
$ws.Range("W4").Value = [double]"1.054585307969091e-10"
$ws.Range("X4").Value = [double]"13.20894894894902"
$ws.Range("Y4").Value = [double]"9.369669669669726"
$ws.Range("Z4").Value = [double]"17.04822822822832"

$ws.Range("F5").Value = [double]"23.55000000000024"
$ws.Range("H5").Value = [double]"0.0002053397409436108"
$ws.Range("I5").Value = [double]"0.0002053397409436108"
$ws.Range("L5").Value = [double]"50.19147722104001"
$ws.Range("M5").Value = "[21.137714626249576, 79.24523981583044]"
$ws.Range("N5").Value = [double]"0.001127577451902839"
$ws.Range("O5").Value = [double]"0.001127577451902839"
$ws.Range("P5").Value = [double]"2.761079429281119"
$ws.Range("Q5").Value = "[2.144710900648888, 3.3774479579133505]"
$ws.Range("R5").Value = [double]"1.178279696034679e-11"
$ws.Range("S5").Value = [double]"1.178279696034679e-11"
$ws.Range("T5").Value = [double]"56.99872521809787"
$ws.Range("U5").Value = "[41.17273441942574, 72.82471601677001]"
$ws.Range("V5").Value = [double]"4.287203037023346e-09"
$ws.Range("W5").Value = [double]"4.287203037023346e-09"
$ws.Range("X5").Value = [double]"13.20120120120134"
$ws.Range("Y5").Value = [double]"10.89099099099111"
$ws.Range("Z5").Value = [double]"15.51141141141157"

$ws.Range("F6").Value = [double]"23.55000000000024"
$ws.Range("H6").Value = [double]"0.0004924609353286202"
$ws.Range("I6").Value = [double]"0.0004924609353286202"
$ws.Range("L6").Value = [double]"37.71261761473687"
$ws.Range("M6").Value = "[13.56270234285794, 61.8625328866158]"
$ws.Range("N6").Value = [double]"0.002938964356721474"
$ws.Range("O6").Value = [double]"0.002938964356721474"
$ws.Range("P6").Value = [double]"2.72334258058935"
$ws.Range("Q6").Value = "[2.0566582537014257, 3.3900269074772735]"
$ws.Range("R6").Value = [double]"1.612396882677558e-10"
$ws.Range("S6").Value = [double]"1.612396882677558e-10"
$ws.Range("T6").Value = [double]"57.83169080265269"
$ws.Range("U6").Value = "[45.14726340649075, 70.51611819881462]"
$ws.Range("V6").Value = [double]"7.005507285384738e-12"
$ws.Range("W6").Value = [double]"7.005507285384738e-12"
$ws.Range("X6").Value = [double]"13.34264264264278"
$ws.Range("Y6").Value = [double]"10.84384384384396"
$ws.Range("Z6").Value = [double]"15.84144144144161"

$ws.Range("F7").Value = [double]"23.55000000000024"
$ws.Range("H7").Value = [double]"0.0001829868893320485"
$ws.Range("I7").Value = [double]"0.0001829868893320485"
$ws.Range("L7").Value = [double]"35.0439421177985"
$ws.Range("M7").Value = "[16.07168734306179, 54.01619689253521]"
$ws.Range("N7").Value = [double]"0.0005502376800994302"
$ws.Range("O7").Value = [double]"0.0005502376800994302"
$ws.Range("P7").Value = [double]"2.962342622303888"
$ws.Range("Q7").Value = "[2.333395144107734, 3.591290100500043]"
$ws.Range("R7").Value = [double]"2.641220575583247e-12"
$ws.Range("S7").Value = [double]"2.641220575583247e-12"
$ws.Range("T7").Value = [double]"48.45918055602299"
$ws.Range("U7").Value = "[37.469189168884164, 59.44917194316182]"
$ws.Range("V7").Value = [double]"1.867372922959021e-11"
$ws.Range("W7").Value = [double]"1.867372922959021e-11"
$ws.Range("X7").Value = [double]"12.44684684684698"
$ws.Range("Y7").Value = [double]"10.0894894894896"
$ws.Range("Z7").Value = [double]"14.80420420420436"
